$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 662.05884
$ws.Range("I17").Value = 641
$ws.Range("J17").Value = 692.1429000000001
$ws.Range("K17").Value = 1923
$ws.Range("L17").Value = 2076.4287
$ws.Range("M17").Value = -1755
$ws.Range("N17").Value = -2412.4287

$ws.Range("H62").Value = 6780.2
$ws.Range("I62").Value = 5451
$ws.Range("J62").Value = 7263.5454
$ws.Range("K62").Value = 5451
$ws.Range("L62").Value = 7263.5454
$ws.Range("M62").Value = -4827
$ws.Range("N62").Value = -8511.545399999999

$ws.Range("H65").Value = 6780.2
$ws.Range("I65").Value = 5451
$ws.Range("J65").Value = 7263.5454
$ws.Range("K65").Value = 27255
$ws.Range("L65").Value = 36317.727
$ws.Range("M65").Value = -24135
$ws.Range("N65").Value = -42557.727

$ws.Range("H115").Value = 575.125
$ws.Range("I115").Value = 740.2
$ws.Range("K115").Value = 2220.6
$ws.Range("M115").Value = -653.6000000000004

$ws.Range("H127").Value = 5706.5
$ws.Range("I127").Value = 5706.5
$ws.Range("K127").Value = 17119.5
$ws.Range("M127").Value = -12159.5

$ws.Range("H131").Value = 8332.666999999999
$ws.Range("I131").Value = 4997
$ws.Range("J131").Value = 9285.714
$ws.Range("K131").Value = 14991
$ws.Range("L131").Value = 27857.142
$ws.Range("M131").Value = -9951
$ws.Range("N131").Value = -37937.142

$ws.Range("H137").Value = 3771.2856
$ws.Range("I137").Value = 4073.5833
$ws.Range("J137").Value = 3544.5625
$ws.Range("K137").Value = 12220.7499
$ws.Range("L137").Value = 10633.6875
$ws.Range("M137").Value = -9670.749899999999
$ws.Range("N137").Value = -15733.6875

$ws.Range("H138").Value = 4774.2666
$ws.Range("I138").Value = 2997
$ws.Range("J138").Value = 5047.6924
$ws.Range("K138").Value = 8991
$ws.Range("L138").Value = 15143.0772
$ws.Range("M138").Value = -3851
$ws.Range("N138").Value = -25423.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 185.75
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 207.2
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 207.2
$ws.Range("M4").Value = -34
$ws.Range("N4").Value = -439.2

$ws.Range("H6").Value = 1538889.9
$ws.Range("I6").Value = 1538889.9
$ws.Range("K6").Value = 1538889.9
$ws.Range("M6").Value = -1538716.9

$ws.Range("H8").Value = 5250
$ws.Range("I8").Value = 2500
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -2356
$ws.Range("N8").Value = -8288

$ws.Range("H32").Value = 3198.8635
$ws.Range("I32").Value = 2652.439
$ws.Range("K32").Value = 2652.439
$ws.Range("M32").Value = -2365.439

$ws.Range("H37").Value = 7006.8
$ws.Range("I37").Value = 7006.8
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 7006.8
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -6733.8
$ws.Range("N37").ClearContents()

$ws.Range("H44").Value = 24999
$ws.Range("J44").Value = 24999
$ws.Range("L44").Value = 24999
$ws.Range("N44").Value = -25975

$ws.Range("H55").Value = 24999
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25629

$ws.Range("H80").Value = 99999.664
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 99999.664
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 99999.664
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -101995.664

$ws.Range("H83").Value = 99999.664
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 99999.664
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 299998.992
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -309982.992

$ws.Range("H88").Value = 3500
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 3500
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H122").Value = 2859.875
$ws.Range("I122").Value = 1649.1875
$ws.Range("J122").Value = 5281.25
$ws.Range("K122").Value = 4947.5625
$ws.Range("L122").Value = 15843.75
$ws.Range("M122").Value = -2497.5625
$ws.Range("N122").Value = -20743.75

$ws.Range("H132").Value = 12893.617
$ws.Range("I132").Value = 11824.469
$ws.Range("K132").Value = 35473.407
$ws.Range("M132").Value = -32943.407

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 240.2
$ws.Range("I22").Value = 240.2
$ws.Range("K22").Value = 240.2
$ws.Range("M22").Value = -67.19999999999999

$ws.Range("H86").Value = 3493.75
$ws.Range("I86").Value = 3278.5715
$ws.Range("K86").Value = 3278.5715
$ws.Range("M86").Value = -2155.5715

$ws.Range("H89").Value = 3493.75
$ws.Range("I89").Value = 3278.5715
$ws.Range("K89").Value = 16392.8575
$ws.Range("M89").Value = -10776.8575

$ws.Range("H105").Value = 35808.547
$ws.Range("I105").Value = 48414.145
$ws.Range("J105").Value = 13748.75
$ws.Range("K105").Value = 48414.145
$ws.Range("L105").Value = 13748.75
$ws.Range("M105").Value = -46667.145
$ws.Range("N105").Value = -17242.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 37178.91
$ws.Range("I86").Value = 19375
$ws.Range("J86").Value = 47352.57
$ws.Range("K86").Value = 19375
$ws.Range("L86").Value = 47352.57
$ws.Range("M86").Value = -18252
$ws.Range("N86").Value = -49598.57

$ws.Range("H89").Value = 37178.91
$ws.Range("I89").Value = 19375
$ws.Range("J89").Value = 47352.57
$ws.Range("K89").Value = 96875
$ws.Range("L89").Value = 236762.85
$ws.Range("M89").Value = -91259
$ws.Range("N89").Value = -247994.85

$ws.Range("H132").Value = 2894.7778
$ws.Range("I132").Value = 2894.7778
$ws.Range("K132").Value = 8684.3334
$ws.Range("M132").Value = -6154.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1041902.75
$ws.Range("I2").Value = 1736149.5
$ws.Range("J2").Value = 532.6667
$ws.Range("K2").Value = 10416897
$ws.Range("L2").Value = 3196.0002
$ws.Range("M2").Value = -10416784
$ws.Range("N2").Value = -3422.0002

$ws.Range("H7").Value = 7693292
$ws.Range("I7").Value = 14286037
$ws.Range("K7").Value = 42858111
$ws.Range("M7").Value = -42857999

$ws.Range("H32").Value = 34849132
$ws.Range("J32").Value = 34849132
$ws.Range("L32").Value = 104547396
$ws.Range("N32").Value = -104547962

$ws.Range("H46").Value = 1111707.5
$ws.Range("I46").Value = 734
$ws.Range("J46").Value = 2500424.2
$ws.Range("K46").Value = 2202
$ws.Range("L46").Value = 7501272.600000001
$ws.Range("M46").Value = -2111
$ws.Range("N46").Value = -7501454.600000001

$ws.Range("H113").Value = 881.5
$ws.Range("I113").Value = 881.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2644.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -474.5
$ws.Range("N113").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws.Range("H140").Value = 1987.5714
$ws.Range("I140").Value = 1802.6
$ws.Range("K140").Value = 5407.799999999999
$ws.Range("M140").Value = -227.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 96.89286
$ws.Range("I2").Value = 42.75
$ws.Range("J2").Value = 137.5
$ws.Range("K2").Value = 42.75
$ws.Range("L2").Value = 137.5
$ws.Range("M2").Value = 70.25
$ws.Range("N2").Value = -363.5

$ws.Range("H7").Value = 7297.3335
$ws.Range("J7").Value = 7297.3335
$ws.Range("L7").Value = 7297.3335
$ws.Range("N7").Value = -7521.3335

$ws.Range("H8").Value = 7297.3335
$ws.Range("J8").Value = 7297.3335
$ws.Range("L8").Value = 7297.3335
$ws.Range("N8").Value = -7575.3335

$ws.Range("H100").Value = 28055
$ws.Range("J100").Value = 28055
$ws.Range("L100").Value = 28055
$ws.Range("N100").Value = -30219

$ws.Range("H122").Value = 503997
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 42500
$ws.Range("I122").Value = 39998.75
$ws.Range("J122").Value = 47502.5
$ws.Range("K122").Value = 119996.25
$ws.Range("L122").Value = 142507.5
$ws.Range("M122").Value = -117546.25
$ws.Range("N122").Value = -147407.5

$ws.Range("H140").Value = 20429
$ws.Range("J140").Value = 20429
$ws.Range("L140").Value = 20429
$ws.Range("N140").Value = -30789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2525
$ws.Range("I81").Value = 2525
$ws.Range("K81").Value = 5050
$ws.Range("M81").Value = -3989

$ws.Range("H84").Value = 2525
$ws.Range("I84").Value = 2525
$ws.Range("K84").Value = 25250
$ws.Range("M84").Value = -19946

$ws.Range("H104").Value = 27726.428
$ws.Range("J104").Value = 27726.428
$ws.Range("L104").Value = 27726.428
$ws.Range("N104").Value = -34714.428

$ws.Range("H118").Value = 39994.332
$ws.Range("I118").Value = 39989
$ws.Range("J118").Value = 39997
$ws.Range("K118").Value = 39989
$ws.Range("L118").Value = 39997
$ws.Range("M118").Value = -38332
$ws.Range("N118").Value = -43311
